$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2340425531914894
$ws.Range("C2").Value = 0.5531914893617021
$ws.Range("P2").Value = 0.1276595744680851
$ws.Range("S2").Value = 0.0851063829787234
$ws.Range("P3").Value = 0.8846153846153846
$ws.Range("S3").Value = 0.1153846153846154
$ws.Range("P4").Value = 1
$ws.Range("B6").Value = 0.02173913043478261
$ws.Range("F6").Value = 0.04347826086956522
$ws.Range("J6").Value = 0.3478260869565217
$ws.Range("Q6").Value = 0.1304347826086956
$ws.Range("R6").Value = 0.02173913043478261
$ws.Range("S6").Value = 0.4347826086956522
$ws.Range("B7").Value = 0.01886792452830189
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.05660377358490566
$ws.Range("J7").Value = 0.07547169811320754
$ws.Range("O7").Value = 0.01886792452830189
$ws.Range("Q7").Value = 0.1132075471698113
$ws.Range("R7").Value = 0.09433962264150944
$ws.Range("S7").Value = 0.6037735849056604
$ws.Range("B8").Value = 0.05521472392638037
$ws.Range("D8").Value = 0.006134969325153374
$ws.Range("E8").Value = 0.006134969325153374
$ws.Range("F8").Value = 0.05521472392638037
$ws.Range("J8").Value = 0.06748466257668712
$ws.Range("O8").Value = 0.03067484662576687
$ws.Range("Q8").Value = 0.147239263803681
$ws.Range("R8").Value = 0.049079754601227
$ws.Range("S8").Value = 0.5828220858895705
$ws.Range("B9").Value = 0.06060606060606061
$ws.Range("F9").Value = 0.07575757575757576
$ws.Range("J9").Value = 0.07575757575757576
$ws.Range("Q9").Value = 0.196969696969697
$ws.Range("R9").Value = 0.06060606060606061
$ws.Range("S9").Value = 0.5303030303030303
$ws.Range("B10").Value = 0.09049773755656108
$ws.Range("D10").Value = 0.02262443438914027
$ws.Range("F10").Value = 0.04524886877828054
$ws.Range("J10").Value = 0.07239819004524888
$ws.Range("O10").Value = 0.004524886877828055
$ws.Range("Q10").Value = 0.1855203619909502
$ws.Range("R10").Value = 0.07239819004524888
$ws.Range("S10").Value = 0.5067873303167421
$ws.Range("G11").Value = 0.1733333333333333
$ws.Range("J11").Value = 0.09333333333333334
$ws.Range("K11").Value = 0.24
$ws.Range("L11").Value = 0.4266666666666667
$ws.Range("S11").Value = 0.06666666666666667
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.25
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.02222222222222222
$ws.Range("H15").Value = 0.06666666666666667
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.3555555555555556
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.02222222222222222
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.2222222222222222
$ws.Range("I16").Value = 0.1944444444444444
$ws.Range("J16").Value = 0.1944444444444444
$ws.Range("K16").Value = 0.05555555555555555
$ws.Range("M16").Value = 0.05555555555555555
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.2222222222222222
$ws.Range("F17").Value = 0.01123595505617977
$ws.Range("H17").Value = 0.1797752808988764
$ws.Range("I17").Value = 0.1460674157303371
$ws.Range("J17").Value = 0.2808988764044944
$ws.Range("K17").Value = 0.1235955056179775
$ws.Range("M17").Value = 0.02247191011235955
$ws.Range("O17").Value = 0.03370786516853932
$ws.Range("S17").Value = 0.2022471910112359
$ws.Range("H18").Value = 0.2058823529411765
$ws.Range("I18").Value = 0.1176470588235294
$ws.Range("J18").Value = 0.2058823529411765
$ws.Range("K18").Value = 0.08823529411764706
$ws.Range("O18").Value = 0.1470588235294118
$ws.Range("S18").Value = 0.2352941176470588
$ws.Range("F19").Value = 0.02682926829268293
$ws.Range("H19").Value = 0.3170731707317073
$ws.Range("I19").Value = 0.0975609756097561
$ws.Range("J19").Value = 0.2439024390243902
$ws.Range("K19").Value = 0.08780487804878048
$ws.Range("M19").Value = 0.04146341463414634
$ws.Range("N19").Value = 0.002439024390243902
$ws.Range("O19").Value = 0.04146341463414634
$ws.Range("S19").Value = 0.1414634146341463
